$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# C4: estimated-completion date moves up a day (43788 -> 43787)
$ws.Range("C4").Value = 43787

# A new persona's daily entries push the open rows (8-10) to carry real
# dates now; copy the date formatting from an existing date cell (A7) so
# the cell style matches the template instead of creating a new style.
$ws.Range("A7").Copy()
$ws.Range("A8:A10").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

$ws.Range("A8").Value = 43785
$ws.Range("A9").Value = 43786
$ws.Range("A10").Value = 43787

# That persona logged 1 point on 43786 (row 9)
$ws.Range("B9").Value = 1

# Row 11 no longer holds the trailing date/formulas - clear it back out
$ws.Range("A11").ClearContents()
$ws.Range("F11").ClearContents()
$ws.Range("G11").ClearContents()

# Restore the view: scrolled back to the top, selection on C5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C5").Select()
